$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.939.50'
$ws.Cells.Item(2, 5).Value = '  +0.83%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.631.28'
$ws.Cells.Item(3, 5).Value = '  +1.79%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.28%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '214.33'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.88%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.11%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.29%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '28.52'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.57%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.258'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.59%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0608'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.85%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.11%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.862.82'
$ws.Cells.Item(12, 5).Value = '  +1.66%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.624.58'
$ws.Cells.Item(13, 5).Value = '  +1.40%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.562'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.48%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Chainlink'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '9.16'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +13.57%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '29.952.08'
$ws.Cells.Item(16, 5).Value = '  +0.82%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.84'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.46%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '64.08'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.07%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '241.86'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.10%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  +0.29%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.25%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.13'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.39%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.75'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.70%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '2.17'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.96%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '158.43'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.83%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '15.50'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.22%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  +0.58%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.59'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +1.75%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.20%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0487'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.46%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.15%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.37'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.05%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '3.17'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.23%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.425.67'
$ws.Cells.Item(34, 5).Value = '  +0.10%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.65'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +4.69%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -1.50%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.81%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.27%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.14%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '75.70'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +12.30%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.551'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.76%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.00'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.67%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.826'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.73%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0489'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -1.27%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.32%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.49%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '52.09'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -5.68%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'FraxShare'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '5.34'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.67%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(49, 4).Value = '1.769.57'
$ws.Cells.Item(49, 5).Value = '  +1.70%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +11.07%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '90.05'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +3.86%  '
